$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H2 from numeric 1 to text "1;2;4"
$ws.Range("H2").Value = "1;2;4"

# Add E3 = "POSL"
$ws.Range("E3").Value = "POSL"

# Add F4 = "OPEN"
$ws.Range("F4").Value = "OPEN"

# Update selection to H2 (active cell H2, selected range H2)
$ws.Range("H2").Select()
